$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.255.81"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.326.92"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'301.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'100.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").Value = "'0.509"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("D10").Value = "'36.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.32%  "
$ws.Range("D11").Value = "'0.0793"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "'17.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "'6.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").Value = "2.685.32"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "2.337.97"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "43.177.99"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'12.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'6.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "'68.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "'236.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "'2.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.70%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'25.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("D28").Value = "'168.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "'34.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").Value = "'9.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -10.49%  "
$ws.Range("D32").Value = "'5.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.27%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "'17.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'4.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").Value = "'2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'0.0698"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").Value = "'0.103"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").Value = "'1.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.001.73"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0292"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("E44").Value = "  -4.61%  "
$ws.Range("D45").Value = "'10.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "'17.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "'2.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("D48").Value = "'55.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").Value = "'1.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("D50").Value = "2.551.73"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'71.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.99%  "
